$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Milena Raquel Pereira
$ws.Cells.Item(13, 1).Value = 13
$ws.Cells.Item(13, 2).Value = "Milena Raquel Pereira"
$ws.Cells.Item(13, 3).Value = "milenaraquelpereira@acquire.com.br"
$ws.Cells.Item(13, 4).Value = "27/07/1965"

# Row 14: Sérgio Diego Da Rocha
# "06/01/2007" parses as a valid day/month date, so Excel would otherwise
# auto-convert it to a date serial. Force it to stay as literal text, then
# restore the cell's style to the same (default) style used by the other
# data rows so no extra formatting is introduced.
$ws.Cells.Item(14, 1).Value = 14
$ws.Cells.Item(14, 2).Value = "Sérgio Diego Da Rocha"
$ws.Cells.Item(14, 3).Value = "sergiodiegodarocha@catsfeelings.com.br"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "06/01/2007"
$ws.Cells.Item(14, 4).Style = $ws.Cells.Item(2, 4).Style
